$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# --- Row 55 ---------------------------------------------------------------
# Copy formatting from the last existing data row (54) into the new cells,
# then overwrite their values (mirrors typing directly below existing rows,
# which is how Excel inherits the row-5 "data" style (s="5")).
$ws.Range("B54").Copy()
$ws.Range("B55").PasteSpecial(-4122)
$ws.Range("B54").Copy()
$ws.Range("B56").PasteSpecial(-4122)

$ws.Range("B55").Value = "fm51"
$ws.Range("B56").Value = "fm52"

$ws.Range("C54").Copy()
$ws.Range("C55").PasteSpecial(-4122)
$ws.Range("C54").Copy()
$ws.Range("C56").PasteSpecial(-4122)

$ws.Range("C55").Value = "Zero intermediate loss example (max deductible) #54"
$ws.Range("C56").Value = "Max deductible not being applied for some samples #55"

$ws.Range("D54").Copy()
$ws.Range("D55").PasteSpecial(-4122)
$ws.Range("D55").Value = "All"

$ws.Range("H54").Copy()
$ws.Range("H55").PasteSpecial(-4122)
$ws.Range("H54").Copy()
$ws.Range("H56").PasteSpecial(-4122)
$ws.Range("I54").Copy()
$ws.Range("I55").PasteSpecial(-4122)
$ws.Range("I54").Copy()
$ws.Range("I56").PasteSpecial(-4122)

$ws.Range("H55").Value = "in progress"
$ws.Range("I55").Value = "in progress"
$ws.Range("H56").Value = "in progress"
$ws.Range("I56").Value = "in progress"

$excel.CutCopyMode = $false

# --- View state: select the newly added H56:I56 range, matching the author
#     ending their edit session with that range selected -------------------
$ws.Range("H56:I56").Select()
